$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: ip_address_list
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ip_address_list")

# Row 1 - truncate the note in D1, keep trailing space
$ws1.Range("D1").Value = "FortiClient Austin: "

# Row 3 - truncate the note in D3 down to the first line
$ws1.Range("D3").Value = "PC:`t10.96.205.175"

# Row 5 - truncate the note in D5 down to the first line
$ws1.Range("D5").Value = "PC:192.168.14.240"

# Row 9 - add a new note in D9
$ws1.Range("D9").Value = "saggggggggggf"

# Row 11 - shorten the name/ip/note by one trailing character
$ws1.Range("A11").Value = "527_Tei"
$ws1.Range("B11").Value = "10.101.28.17"
$ws1.Range("D11").Value = "XG-X2900:`t`t10.101.28.175`nOP:`t`t10.101.28."
$ws1.Rows.Item(11).AutoFit()

# ---------------------------------------------------------------
# Sheet 2: ip_adress_fav_list
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")

# Row 1 - add new note, mark favourite flag as boolean TRUE
$ws2.Range("D1").Value = "saggggggggggf"
$ws2.Range("E1").Value = $true

# Row 2 - replace the 514_Teleflex entry with the updated 527_Tei entry
$ws2.Range("A2").Value = "527_Tei"
$ws2.Range("B2").Value = "10.101.28.17"
$ws2.Range("D2").Value = "XG-X2900:`t`t10.101.28.175`nOP:`t`t10.101.28."
$ws2.Range("E2").Value = $true
$ws2.Rows.Item(2).AutoFit()

# Row 3 - no longer needed, remove entirely
$ws2.Rows.Item(3).Delete()

# ---------------------------------------------------------------
# Sheet 3: disk_list
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("disk_list")

$ws3.Range("F3").Value = "f"
$ws3.Range("F4").Value = "Druha s" + [char]0xED + [char]0x9D + "t, ixonah"
$ws3.Range("F5").Value = "Druha s" + [char]0xED + [char]0x9D + "t, ixonah"
